$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$t17 = @"
Ensinar a linguagem gráfica normalizada internacionalmente para representação de máquinas e equipamentos que integram os processos de engenharia.Desenvolver o raciocínio espacial e a criatividade de representação.
"@

$t20 = @"
5840820 - Gustavo Aristides Santana Martinez
"@

$t21 = @"
Programa resumido:
"@

$t22 = @"
-Introdução-Teoria Elementar do Desenho Projetivo-Projeções Ortogonais pelo 1º Diedro-Projeções Ortogonais pelo 3º Diedro-Leitura e Interpretação de Desenhos-Escalas-Desenhos com Instrumentos-Cortes e Representações Convencionais-Projeções Auxiliares-Cotação-Desenhos de Conjuntos e Detalhes-Aplicação de Tolerâncias e Ajustes -Símbolos de Acabamento Superficial-Desenho de Elementos de Máquina-Desenho de Equipamentos e Acessórios
"@

$t23 = @"
Short syllabus:
"@

$t24 = @"
Programa:
"@

$t25 = @"
1 - INTRODUÇÃOApresentação e definição da disciplina, destacando a importância do desenho na engenharia; Normas ABNT e ISO.2 - TEORIA ELEMENTAR DO DESENHO PROJETIVORepresentação de vistas como sistema internacional; representação de arestas visíveis e invisíveis; linhas de centro e eixos de simetria.3 - PROJEÇÕES ORTOGONAIS PELO 1º DIEDROPrincípio fundamental; projeções principais; rebatimentos convencionados.4 - PROJEÇÕES ORTOGONAIS PELO 3º DIEDROPrincípio fundamental; projeções principais; rebatimentos convencionados.5 - LEITURA E INTERPRETAÇÃO DE DESENHOSLeitura por meio de esboço em perspectiva e mediante construção de modelos.6 - ESCALASDefinição e normalização7 - DESENHOS COM INSTRUMENTOSRegras para emprego dos esquadros, compasso e régua "T"; disposição do desenho nas folhas padronizadas.8 - CORTES E REPRESENTAÇÕES CONVENCIONAISPrincípios fundamentais; aplicações; tipos normalizados; representações e regras para traçado; seções e rupturas.9 - PROJEÇÕES AUXILIARESPrincípios fundamentais; finalidades e aplicações; representações normalizadas.10 - COTAÇÃORegras de colocação e distribuição de cotas.11 - DESENHOS DE CONJUNTOS E DETALHESDefinições; tipos recomendados de legenda e lista de peça; formas de numeração de desenhos; regras práticas para execução e verificação de desenhos.12 - APLICAÇÃO DE TOLERÂNCIAS E AJUSTESDefinição e finalidades; sistema ISO; uso de tabelas e indicação nos desenhos.13 - SÍMBOLOS DE ACABAMENTO SUPERFICIALDefinição; simbologia normalizada; aplicações.14 - DESENHO DE ELEMENTOS DE MÁQUINADefinições, aplicações, tipos, proporções e representações convencionais de: roscas, parafusos, porcas, arruelas, polias, correias e chavetas.15 - DESENHO DE EQUIPAMENTOS E ACESSÓRIOSDesenho de conjunto e detalhes envolvendo elementos de ligação e de máquinas com aplicação de tabelas e catálogos.
"@

$t26 = @"
Syllabus:
"@

$t27 = @"
Avaliação:
"@

$t28 = @"
Método:
"@

$t29 = @"
A avaliação é continuada e constará de duas provas objetivas (Pi) realizadas ao longo do curso (antes da recuperação), bem como de exercícios práticos realizados em sala de aula e extra classe (TC/TS).
"@

$t30 = @"
Critério:
"@

$t31 = @"
NOTA FINAL = [(MédiaTC/TS)x0,2] + [(MédiaPi)x0,8]
"@

$t32 = @"
Norma de recuperação:
"@

$t33 = @"
- A recuperação deverá consistir de uma prova englobando a matéria toda do semestre.- A média final (pós-recuperação) deverá ser composta por uma média simples entre a nota do semestre (nota final) e a da prova de recuperação.
"@

$t34 = @"
Bibliografia:
"@

$t35 = @"
1 - ABNT - COLETÂNEA DE NORMAS DE DESENHO TÉCNICONormas Técnicas publicadas pela ABNT2 - DESENHO BÁSICO NA ENGENHARIARibeiro, Antonio Clélio - Apostila publicada pela FAENQUIL3 - FUNDAMENTOS DE DIBUJO EM INGENIERIALuzader, Warren J. - Ed. Comp. Editorial Continental - México4 - MANUAL DE DESENHO TÉCNICOManfé, G./ Scarato, G./ Pozza, R. - Ed. Renovada Livros Culturais Ltda.5 - EXPRESSÃO GRÁFICA - DESENHO TÉCNICOHoelsher, R. P./ Springer, C.H./ Dobrovolny, J.S. - Ed. LTC Editora S.A.6 - DESENHO TÉCNICOFrench, Thomas E. - Editora Globo7 - DESENHO TÉCNICOBachmann, A./ Forberg, R - Editora Globo8 - DESENHISTA DE MÁQUINASEscola PRO-TEC
"@

# --- Row 10: Objetivos value (B10/C10) ---
$ws.Range("B10").Value = $t17
$ws.Range("C10").Value = $t17

# --- Row 13: remove A13 label (now part of row 12 "Docentes responsaveis" block);
#     B13/C13 keep teacher name, row height reverts to default ---
$ws.Range("A13").Clear()
$ws.Range("B13").Value = $t20
$ws.Range("C13").Value = $t20
$ws.Rows.Item(13).AutoFit()

# --- Row 14: A14 becomes "Programa resumido:" label; B14/C14 get the new summary text ---
# (B14 is a brand-new cell; copy column-B formatting first so it gets style index 2, not A's style)
$ws.Range("B10").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A14").Value = $t21
$ws.Range("B14").Value = $t22
$ws.Range("C14").Value = $t22
$ws.Rows.Item(14).RowHeight = 60

# --- Row 15: A15 becomes "Short syllabus:"; clear stray B15/C15 data, height -> 60 ---
$ws.Range("A15").Value = $t23
$ws.Range("B15").Clear()
$ws.Range("C15").Clear()
$ws.Rows.Item(15).RowHeight = 60

# --- Row 16: A16 becomes "Programa:"; B16/C16 get the new full program text, height -> 120 ---
# (B16 is a brand-new cell; copy column-B formatting first so it gets style index 2, not A's style)
$ws.Range("B10").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A16").Value = $t24
$ws.Range("B16").Value = $t25
$ws.Range("C16").Value = $t25
$ws.Rows.Item(16).RowHeight = 120

# --- Row 17: A17 becomes "Syllabus:", height -> 120 ---
$ws.Range("A17").Value = $t26
$ws.Rows.Item(17).RowHeight = 120

# --- Row 18: A18 becomes "Avaliacao:"; clear stray B18/C18 data, height reverts to default ---
$ws.Range("A18").Value = $t27
$ws.Range("B18").Clear()
$ws.Range("C18").Clear()
$ws.Rows.Item(18).AutoFit()

# --- Row 19: A19 "Metodo:"; B19/C19 get avaliacao continuada text ---
$ws.Range("A19").Value = $t28
$ws.Range("B19").Value = $t29
$ws.Range("C19").Value = $t29

# --- Row 20: A20 "Criterio:"; B20/C20 get NOTA FINAL text ---
$ws.Range("A20").Value = $t30
$ws.Range("B20").Value = $t31
$ws.Range("C20").Value = $t31

# --- Row 21: A21 "Norma de recuperacao:"; B21/C21 get recuperacao text, height -> 60 ---
$ws.Range("A21").Value = $t32
$ws.Range("B21").Value = $t33
$ws.Range("C21").Value = $t33
$ws.Rows.Item(21).RowHeight = 60

# --- Row 22 (NEW): copy formatting from row 21, then set "Bibliografia:" label and text, height -> 120 ---
$ws.Range("A21:C21").Copy()
$ws.Range("A22:C22").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A22").Value = $t34
$ws.Range("B22").Value = $t35
$ws.Range("C22").Value = $t35
$ws.Rows.Item(22).RowHeight = 120
